# gt_buttons.xlsx - "Add files via upload" edit
# - Reorders the category-3/category-4 button rows (17-25) into the new
#   grouping used by the app (Stop/Traffic light/Roundabout/Speed bump/Yield
#   sign become Row 2; Start of highway/urban city/ALCA become Row 3;
#   End of highway/urban city/ALCA become Row 4).
# - Adds two brand new buttons: "Start of ALCA available" (Row 3) and
#   "End of ALCA available" (Row 4).
# - Updates selection/column width cosmetics to match the refreshed sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for column A (button name) / column C (Row group) for rows 17-27.
$ws.Range("A17").Value2 = "Stop sign"
$ws.Range("C17").Value2 = 2

$ws.Range("A18").Value2 = "Traffic light"
$ws.Range("C18").Value2 = 2

$ws.Range("A19").Value2 = "Roundabout"
$ws.Range("C19").Value2 = 2

$ws.Range("A20").Value2 = "Speed bump"
$ws.Range("C20").Value2 = 2

$ws.Range("A21").Value2 = "Yield sign"
$ws.Range("C21").Value2 = 2

$ws.Range("A22").Value2 = "Start of highway"
$ws.Range("C22").Value2 = 3

$ws.Range("A23").Value2 = "Start of urban/city"
$ws.Range("C23").Value2 = 3

$ws.Range("A24").Value2 = "Start of ALCA available"
$ws.Range("C24").Value2 = 3

$ws.Range("A25").Value2 = "End of highway"
$ws.Range("C25").Value2 = 4

# New rows added at the bottom of the table.
$ws.Range("A26").Value2 = "End of urban/city"
$ws.Range("C26").Value2 = 4

$ws.Range("A27").Value2 = "End of ALCA available"
$ws.Range("C27").Value2 = 4

# Column A is no longer auto best-fit; it now has an explicit custom width.
$ws.Columns.Item(1).ColumnWidth = 20.25

# Selection moved from the old scratch cell (G26) to A2.
$ws.Range("A2").Select()
